# The edit inserts a new data row at sheet row 19 (pushing the existing
# rows 19-82 down to 20-83, same as the "weekly" price log growing by one
# new daily record at the top of this variety block), and fills the new
# row with its own Fecha/Variedad/Calidad/Volumen/Precio values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push existing rows 19..82 down to 20..83, leaving row 19 blank for the
# new record.
$ws.Rows.Item(19).Insert()

# Populate the newly inserted row 19 with the new record. Columns that are
# constant for every row in this sub-block (market/region/product taxonomy,
# unit, province) are carried over just like on every other row.
$ws.Range("A19").Value = 4
$ws.Range("B19").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C19").Value = "Los Lagos"
$ws.Range("D19").Value = 44414
$ws.Range("E19").Value = 10
$ws.Range("F19").Value = "Fruta"
$ws.Range("G19").Value = 100102
$ws.Range("H19").Value = "Cítricos"
$ws.Range("I19").Value = 100102004
$ws.Range("J19").Value = "Mandarina"
$ws.Range("K19").Value = "Clementina"
$ws.Range("L19").Value = "Primera"
$ws.Range("M19").Value = 600
$ws.Range("N19").Value = 6500
$ws.Range("O19").Value = 6500
$ws.Range("P19").Value = 6500
$ws.Range("Q19").Value = "$/bandeja 10 kilos"
$ws.Range("R19").Value = "Provincia de Limarí"
$ws.Range("S19").Value = 650
$ws.Range("T19").Value = 10
